$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: title/link update (Follow up post)
$ws.Range("D9").Value = "10년차 고등학교 수학 강사의 Data Science 도전? – Follow up"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/math-teacher-for-data-science-2/#utm_source=rss&utm_medium=rss&utm_campaign=math-teacher-for-data-science-2"

# Row 45: title/link update (Ensemble model - CatBoost)
$ws.Range("D45").Value = "Ensemble model - CatBoost"
$ws.Range("E45").Value = "https://dive-into-ds.tistory.com/120"

# Row 51: title/link update (meta tag post)
$ws.Range("D51").Value = "[html] meta 태그란 왜 있는 것인가? meta 태그 정리"
$ws.Range("E51").Value = "https://bskyvision.com/entry/html-meta-%ED%83%9C%EA%B7%B8%EB%9E%80-%EC%99%9C-%EC%9E%88%EB%8A%94-%EA%B2%83%EC%9D%B8%EA%B0%80-meta-%ED%83%9C%EA%B7%B8-%EC%A0%95%EB%A6%AC"

# Row 52: title update (숨은 DS)
$ws.Range("D52").Value = "숨은 DS"
